# Scheduled-runner style update of cached "goblin profits" price/profit
# figures across the leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM,
# LTW, WVR). These columns (H..N) hold externally-sourced market prices
# and the profit math derived from them; they carry no formulas, so the
# refreshed numbers are written directly as literal values.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H86").Value = 1775
$ws.Range("I86").Value = 1272.5
$ws.Range("J86").Value = 2205.7144
$ws.Range("K86").Value = 1272.5
$ws.Range("L86").Value = 2205.7144
$ws.Range("M86").Value = -149.5
$ws.Range("N86").Value = -4451.7144

$ws.Range("H89").Value = 1775
$ws.Range("I89").Value = 1272.5
$ws.Range("J89").Value = 2205.7144
$ws.Range("K89").Value = 6362.5
$ws.Range("L89").Value = 11028.572
$ws.Range("M89").Value = -746.5
$ws.Range("N89").Value = -22260.572

$ws.Range("H132").Value = 1803.6957
$ws.Range("I132").Value = 1324.6046
$ws.Range("K132").Value = 3973.8138
$ws.Range("M132").Value = -1443.8138

$ws.Range("H138").Value = 3459.2185
$ws.Range("J138").Value = 3942.9048
$ws.Range("L138").Value = 11828.7144
$ws.Range("N138").Value = -22108.7144

$ws.Range("H141").Value = 6000
$ws.Range("I141").Value = 8000
$ws.Range("K141").Value = 24000
$ws.Range("M141").Value = -18820

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 2299.111
$ws.Range("I74").Value = 2345.5334
$ws.Range("J74").Value = 2067
$ws.Range("K74").Value = 2345.5334
$ws.Range("L74").Value = 2067
$ws.Range("M74").Value = -1471.5334
$ws.Range("N74").Value = -3815

$ws.Range("H77").Value = 2299.111
$ws.Range("I77").Value = 2345.5334
$ws.Range("J77").Value = 2067
$ws.Range("K77").Value = 11727.667
$ws.Range("L77").Value = 10335
$ws.Range("M77").Value = -7359.666999999999
$ws.Range("N77").Value = -19071

$ws.Range("H110").Value = 692.25
$ws.Range("I110").Value = 692.25
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 692.25
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1352.75
$ws.Range("N110").Value = $null

# --- BSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H20").Value = 1062
$ws.Range("I20").Value = 1009
$ws.Range("J20").Value = 1199.8
$ws.Range("K20").Value = 1009
$ws.Range("L20").Value = 1199.8
$ws.Range("M20").Value = -762
$ws.Range("N20").Value = -1693.8

$ws.Range("H105").Value = 2502
$ws.Range("I105").Value = 1671.3334
$ws.Range("J105").Value = 3498.8
$ws.Range("K105").Value = 1671.3334
$ws.Range("L105").Value = 3498.8
$ws.Range("M105").Value = 75.66660000000002
$ws.Range("N105").Value = -6992.8

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 2159.4546
$ws.Range("I16").Value = 1351.8334
$ws.Range("J16").Value = 3128.6
$ws.Range("K16").Value = 1351.8334
$ws.Range("L16").Value = 3128.6
$ws.Range("M16").Value = -1064.8334
$ws.Range("N16").Value = -3702.6

$ws.Range("H31").Value = 2925.1516
$ws.Range("I31").Value = 1614.7142
$ws.Range("K31").Value = 1614.7142
$ws.Range("M31").Value = -1319.7142

$ws.Range("H34").Value = 2925.1516
$ws.Range("I34").Value = 1614.7142
$ws.Range("K34").Value = 1614.7142
$ws.Range("M34").Value = -1412.7142

$ws.Range("H51").Value = 25000

$ws.Range("H61").Value = 25000

$ws.Range("H99").Value = 3840
$ws.Range("I99").Value = 2733.3333
$ws.Range("K99").Value = 2733.3333
$ws.Range("M99").Value = -1235.3333

$ws.Range("H113").Value = 2159.4546
$ws.Range("I113").Value = 1351.8334
$ws.Range("J113").Value = 3128.6
$ws.Range("K113").Value = 1351.8334
$ws.Range("L113").Value = 3128.6
$ws.Range("M113").Value = 818.1666
$ws.Range("N113").Value = -7468.6

$ws.Range("H126").Value = 3840
$ws.Range("I126").Value = 2733.3333
$ws.Range("K126").Value = 8199.999899999999
$ws.Range("M126").Value = -5729.999899999999

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H33").Value = 215.11111
$ws.Range("J33").Value = 24
$ws.Range("L33").Value = 144
$ws.Range("N33").Value = -710

$ws.Range("H107").Value = 2300
$ws.Range("J107").Value = 1983.1538
$ws.Range("L107").Value = 5949.4614
$ws.Range("N107").Value = -9789.4614

$ws.Range("H113").Value = 876.0909
$ws.Range("I113").Value = 1013.63635
$ws.Range("J113").Value = 807.3182
$ws.Range("K113").Value = 3040.90905
$ws.Range("L113").Value = 2421.9546
$ws.Range("M113").Value = -870.9090500000002
$ws.Range("N113").Value = -6761.9546

$ws.Range("H121").Value = 1006.6667
$ws.Range("J121").Value = 995
$ws.Range("L121").Value = 2985
$ws.Range("N121").Value = -5605

$ws.Range("H131").Value = 3331.6667
$ws.Range("I131").Value = 823.75
$ws.Range("K131").Value = 2471.25
$ws.Range("M131").Value = 2568.75

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H97").Value = 4753.7334
$ws.Range("I97").Value = 1531
$ws.Range("J97").Value = 33758.332
$ws.Range("K97").Value = 1531
$ws.Range("L97").Value = 33758.332
$ws.Range("M97").Value = -1035
$ws.Range("N97").Value = -34750.332

$ws.Range("H134").Value = 74500
$ws.Range("J134").Value = 74500
$ws.Range("L134").Value = 223500
$ws.Range("N134").Value = -228570

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H100").Value = 6791.6
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459

# --- WVR -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 2022.75
$ws.Range("J81").Value = 2367
$ws.Range("L81").Value = 4734
$ws.Range("N81").Value = -6856

$ws.Range("H84").Value = 2022.75
$ws.Range("J84").Value = 2367
$ws.Range("L84").Value = 23670
$ws.Range("N84").Value = -34278

$ws.Range("H88").Value = 40000
$ws.Range("J88").Value = 40000
$ws.Range("L88").Value = 40000
$ws.Range("N88").Value = -40812

$ws.Range("H91").Value = 40000
$ws.Range("J91").Value = 40000
$ws.Range("L91").Value = 40000
$ws.Range("N91").Value = -42808

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

$ws.Range("H100").Value = 2429.5
$ws.Range("I100").Value = 2038.8334
$ws.Range("K100").Value = 4077.6668
$ws.Range("M100").Value = -3536.6668

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null
